$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.738.75'
$ws.Range("E2").Value = '  +1.58%  '

$ws.Range("D3").Value = '1.881.25'
$ws.Range("E3").Value = '  +1.48%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.006'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.39%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '333.21'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.18%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.006'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.44%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4704'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.25%  '

$ws.Range("E8").Value = '  +0.96%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '47.57'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.16%  '

$ws.Range("E10").Value = '  +1.93%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.027'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.58%  '

$ws.Range("E12").Value = '  +3.73%  '

$ws.Range("D13").Value = '1.884.59'
$ws.Range("E13").Value = '  +0.98%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.976'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.15%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.149'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.15%  '

$ws.Range("E16").Value = '  +0.67%  '

$ws.Range("B17").Value = 'TRON'
$ws.Range("C17").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.06717'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.20%  '

$ws.Range("B18").Value = 'Litecoin'
$ws.Range("C18").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '87.10'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.27%  '

$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.00001046'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.76%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.34'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.55%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.006'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.37%  '

$ws.Range("B22").Value = 'WrappedBTC'
$ws.Range("C22").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D22").Value = '27.739.60'
$ws.Range("E22").Value = '  +1.54%  '

$ws.Range("B23").Value = 'Uniswap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.540'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.85%  '

$ws.Range("E24").Value = '  +1.48%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.315'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.28%  '

$ws.Range("D26").Value = '2.107.71'
$ws.Range("E26").Value = '  +1.09%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '160.16'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.05%  '

$ws.Range("E28").Value = '  +1.32%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.108'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.22%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.593'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.57%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '121.91'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.70%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9859'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.24%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09480'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.39%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.456'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.00%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.619'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.48%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.366'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.11%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06142'
$ws.Range("D37").Style = "Normal"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02265'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.66%  '

$ws.Range("E39").Value = '  +1.36%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.139'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.14%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6006'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.36%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '10.31'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.62%  '

$ws.Range("B44").Value = 'Decentraland'
$ws.Range("C44").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5731'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.06%  '

$ws.Range("B45").Value = 'WEMIXTOKEN'
$ws.Range("C45").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.257'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.73%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.19'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.33%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.947'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.65%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.395'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.54%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06899'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.36%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '114.38'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.75%  '

$ws.Range("B51").Value = 'EOS'
$ws.Range("C51").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.072'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.87%  '
